# Fix: option symbol urls_as_strings was typo
#
# The cells on Sheet1 (A1:A5) contain text that looks like URLs/paths, but
# they must remain plain strings rather than being turned into live
# hyperlinks. Excel (via write_xlsx) had mistakenly created hyperlinks for
# these cells; this edit removes those hyperlinks and updates the text of a
# few cells to carry explicit scheme prefixes (mailto:, internal:, external:)
# instead of being bare/ambiguous strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove every hyperlink object attached to the worksheet (this also drops
# the <hyperlinks> element and associated relationships from the XML).
$ws.Hyperlinks.Delete()

# Update the plain-text contents of the affected cells.
$ws.Cells.Item(2, 1).Value2 = "mailto:write_xlsx@example.com"
$ws.Cells.Item(4, 1).Value2 = "internal:Sheet1!A1"
$ws.Cells.Item(5, 1).Value2 = "external:c:\foo.xlsx"
